# licor_log_file_blank.xlsx edit
# - clear out the previously-filled-in example/demo values so the log file
#   is a blank template again
# - reset the "set value" placeholder cells to underscores
# - update header text (remove date / name that were filled in, tweak font run)
# - move active selection

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Clear the example rows (2-5) that had been filled in with sample data ---
$ws.Range("A2").ClearContents()
$ws.Range("B2").ClearContents()
$ws.Range("G2").ClearContents()

$ws.Range("A3").ClearContents()
$ws.Range("B3").ClearContents()
$ws.Range("G3").ClearContents()

$ws.Range("A4").ClearContents()
$ws.Range("B4").ClearContents()
$ws.Range("G4").ClearContents()
$ws.Range("H4").ClearContents()

$ws.Range("A5").ClearContents()
$ws.Range("B5").ClearContents()
$ws.Range("G5").ClearContents()

# --- Reset the chamber-condition "set value" cells back to blank placeholders ---
$ws.Range("C22").Value = "_________"
$ws.Range("G22").Value = "_________"
$ws.Range("C23").Value = "_________"
$ws.Range("G23").Value = "_________"
$ws.Range("G24").Value = "_________"

# C24 previously held a unique larger font (sz 15); align its formatting with
# the other "set value" cells (C22/C23/G22/G24, which use the sz-20 font)
# before writing its blanked-out placeholder text.
$ws.Range("C22").Copy()
$ws.Range("C24").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("C24").Value = "_________"

# --- Update the print header: drop filled-in date/name, tweak the machine row ---
$ws.PageSetup.LeftHeader = "&`"Calibri (Body),Regular`"&24`nDate:"
$ws.PageSetup.CenterHeader = "&`"Calibri (Body),Regular`"&18Machine (circle):`n&24Ozz   Gib    Alb    Stan&`"-,Regular`"&30    &`"Calibri (Body),Regular`"&24Yat"
$ws.PageSetup.RightHeader = "&`"Calibri (Body),Regular`"&20`nLicor Enthusiast: ___________________   "

# --- Move the active cell selection ---
$ws.Range("E25").Select()
